$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4152
$ws1.Range("F15").Value = 3105

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4152
$ws4.Range("F19").Value = 3105
